$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data validation list for the "form of burial" dropdown (C21):
# was "Urne,Sarg" -> now "die Urne,der Sarg"
$ws.Range("C21").Validation.Formula1 = "die Urne,der Sarg"

# Update the selected value in C21 from "Sarg" to "der Sarg"
$ws.Range("C21").Value = "der Sarg"

# Update the number of songs (Liederanzahl) in C27 from 2 to 3
$ws.Range("C27").Value = 3

# Update the active selection/cell to C21 (was C27)
$ws.Range("C21").Select()
